$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C8").Value = 23
$ws.Range("C9").Value = 266

$ws.Range("C28").Value = 10
$ws.Range("E28").Value = "Complete"

$ws.Range("C31").Value = 5
$ws.Range("E31").Value = "Complete"

$ws.Range("B30").Font.Bold = $true
$ws.Range("B30").Font.Italic = $true

$ws.Range("E29").Select()
